$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 79.23077000000001
$ws.Range("I4").Value = 79.23077000000001
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 79.23077000000001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 34.76922999999999
$ws.Range("N4").ClearContents()

# Row 33
$ws.Range("H33").Value = 732.6957
$ws.Range("I33").Value = 595.0769
$ws.Range("J33").Value = 911.6
$ws.Range("K33").Value = 595.0769
$ws.Range("L33").Value = 911.6
$ws.Range("M33").Value = -366.0769

# Row 74
$ws.Range("H74").Value = 1295
$ws.Range("I74").Value = 1290
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 1290
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = -354
$ws.Range("N74").Value = -3172

# Row 76
$ws.Range("H76").Value = 12000
$ws.Range("I76").Value = 12000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 12000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -11685

# Row 77
$ws.Range("H77").Value = 1295
$ws.Range("I77").Value = 1290
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 6450
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = -1770
$ws.Range("N77").Value = -15860

# Row 79
$ws.Range("H79").Value = 12000
$ws.Range("I79").Value = 12000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 12000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -10908

# Row 86
$ws.Range("H86").Value = 2500
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1377

# Row 89
$ws.Range("H89").Value = 2500
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6884

# Row 96
$ws.Range("H96").Value = 849.9
$ws.Range("I96").Value = 15
$ws.Range("J96").Value = 942.6667
$ws.Range("K96").Value = 45
$ws.Range("L96").Value = 2828.0001
$ws.Range("M96").Value = 1328
$ws.Range("N96").Value = -5574.0001

# Row 111
$ws.Range("H111").Value = 2980.1428
$ws.Range("I111").Value = 3183.25
$ws.Range("J111").Value = 2709.3333
$ws.Range("K111").Value = 9549.75
$ws.Range("L111").Value = 8127.999899999999
$ws.Range("M111").Value = -6482.75

# Row 115
$ws.Range("H115").Value = 107
$ws.Range("I115").Value = 107
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 321
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 1246
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value = 19499.75
$ws.Range("I46").Value = 19333
$ws.Range("J46").Value = 20000
$ws.Range("K46").Value = 19333
$ws.Range("L46").Value = 20000
$ws.Range("M46").Value = -19014
$ws.Range("N46").Value = -20638

# Row 51
$ws.Range("H51").Value = 40000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 40000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41512

# Row 110
$ws.Range("H110").Value = 1063.5714
$ws.Range("I110").Value = 958.6667
$ws.Range("J110").Value = 1693
$ws.Range("K110").Value = 958.6667
$ws.Range("L110").Value = 1693
$ws.Range("M110").Value = 1086.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 301.2857
$ws.Range("I5").Value = 401.33334
$ws.Range("J5").Value = 226.25
$ws.Range("K5").Value = 401.33334
$ws.Range("L5").Value = 226.25
$ws.Range("M5").Value = -288.33334
$ws.Range("N5").Value = -452.25

# Row 95
$ws.Range("H95").Value = 12458
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 12458
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 12458
$ws.Range("N95").Value = -17950

# Row 107
$ws.Range("H107").Value = 2366.7144
$ws.Range("I107").Value = 2366.7144
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2366.7144
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -446.7143999999998
$ws.Range("N107").ClearContents()

# Row 109
$ws.Range("H109").Value = 59996.5
$ws.Range("I109").Value = 59994
$ws.Range("J109").Value = 59999
$ws.Range("K109").Value = 59994
$ws.Range("L109").Value = 59999
$ws.Range("M109").Value = -58607

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 1598.5385
$ws.Range("I2").Value = 978.6667
$ws.Range("J2").Value = 2129.8572
$ws.Range("K2").Value = 978.6667
$ws.Range("L2").Value = 2129.8572
$ws.Range("M2").Value = -865.6667
$ws.Range("N2").Value = -2355.8572

# Row 22
$ws.Range("H22").Value = 1995
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1995
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1995
$ws.Range("N22").Value = -2695
$ws.Range("M22").ClearContents()

# Row 31
$ws.Range("H31").Value = 1621.2222
$ws.Range("I31").Value = 1621.2222
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1621.2222
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1326.2222
$ws.Range("N31").ClearContents()

# Row 34
$ws.Range("H34").Value = 1621.2222
$ws.Range("I34").Value = 1621.2222
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1621.2222
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1419.2222
$ws.Range("N34").ClearContents()

# Row 37
$ws.Range("H37").Value = 22350
$ws.Range("I37").Value = 4700
$ws.Range("J37").Value = 40000
$ws.Range("K37").Value = 4700
$ws.Range("L37").Value = 40000
$ws.Range("M37").Value = -4593
$ws.Range("N37").Value = -40214

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 11753.77
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 11753.77
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 35261.31
$ws.Range("N34").Value = -35429.31

# Row 40
$ws.Range("H40").Value = 600
$ws.Range("I40").Value = 200
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -731
$ws.Range("N40").Value = -4138

# Row 51
$ws.Range("H51").Value = 3520
$ws.Range("I51").Value = 1900
$ws.Range("J51").Value = 4600
$ws.Range("K51").Value = 5700
$ws.Range("L51").Value = 13800
$ws.Range("M51").Value = -5240
$ws.Range("N51").Value = -14720

# Row 68
$ws.Range("H68").Value = 999
$ws.Range("I68").Value = 999
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2997
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2186
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 999
$ws.Range("I71").Value = 999
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8991
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4935
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 385.2
$ws.Range("I97").Value = 737.5
$ws.Range("J97").Value = 150.33333
$ws.Range("K97").Value = 737.5
$ws.Range("L97").Value = 150.33333
$ws.Range("M97").Value = -241.5
$ws.Range("N97").Value = -1142.33333

# Row 102
$ws.Range("H102").Value = 715.2222
$ws.Range("I102").Value = 835.1111
$ws.Range("J102").Value = 595.3333
$ws.Range("K102").Value = 835.1111
$ws.Range("L102").Value = 595.3333
$ws.Range("M102").Value = 786.8889
$ws.Range("N102").Value = -3839.3333

# Row 107
$ws.Range("H107").Value = 754
$ws.Range("I107").Value = 180.4
$ws.Range("J107").Value = 1471
$ws.Range("K107").Value = 180.4
$ws.Range("L107").Value = 1471
$ws.Range("M107").Value = 1739.6
$ws.Range("N107").Value = -5311

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1851
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1851
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1851
$ws.Range("N22").Value = -2441
$ws.Range("M22").ClearContents()

# Row 26
$ws.Range("H26").Value = 4699.8
$ws.Range("I26").Value = 5166.6665
$ws.Range("J26").Value = 3999.5
$ws.Range("K26").Value = 5166.6665
$ws.Range("L26").Value = 3999.5
$ws.Range("M26").Value = -4871.6665
$ws.Range("N26").Value = -4589.5

# Row 27
$ws.Range("H27").Value = 1851
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1851
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1851
$ws.Range("N27").Value = -2065
$ws.Range("M27").ClearContents()

# Row 31
$ws.Range("H31").Value = 10020.167
$ws.Range("I31").Value = 552.5
$ws.Range("J31").Value = 14754
$ws.Range("K31").Value = 552.5
$ws.Range("L31").Value = 14754
$ws.Range("M31").Value = -304.5
$ws.Range("N31").Value = -15250

# Row 82
$ws.Range("H82").Value = 3800
$ws.Range("I82").Value = 3800
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3800
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3439

# Row 85
$ws.Range("H85").Value = 3800
$ws.Range("I85").Value = 3800
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3800
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2552

# Row 97
$ws.Range("H97").Value = 9750
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 9750
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 9750
$ws.Range("N97").Value = -11732

# Row 100
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2000
$ws.Range("N100").Value = -3082

# Row 136
$ws.Range("H136").Value = 556497.7
$ws.Range("I136").Value = 715135.5600000001
$ws.Range("J136").Value = 1265
$ws.Range("K136").Value = 2145406.68
$ws.Range("L136").Value = 3795
$ws.Range("M136").Value = -2142856.68
$ws.Range("N136").Value = -8895

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Row 63
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51248

# Row 66
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156240

# Row 107
$ws.Range("H107").Value = 692.0769
$ws.Range("I107").Value = 582.9167
$ws.Range("J107").Value = 2002
$ws.Range("K107").Value = 1748.7501
$ws.Range("L107").Value = 6006
$ws.Range("M107").Value = 171.2499

# Row 132
$ws.Range("H132").Value = 1533.125
$ws.Range("I132").Value = 1322.8572
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 3968.5716
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -1438.5716
$ws.Range("N132").Value = -14075

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
